$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.8
$ws.Range("K2").Value = 2.63

# Row 3
$ws.Range("G3").Value = 2.25
$ws.Range("I3").Value = 3.75
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 1.91
$ws.Range("L3").Value = 4.33
$ws.Range("X3").Value = 9
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 23
$ws.Range("AH3").Value = 8
$ws.Range("AI3").Value = 17
$ws.Range("AQ3").Value = 51
$ws.Range("AX3").Value = 21
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 351

# Row 4
$ws.Range("R4").Value = 1.48

# Row 6
$ws.Range("R6").Value = 1.36

# Row 7
$ws.Range("G7").Value = 1.8
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 4.75
$ws.Range("L7").Value = 5.5
$ws.Range("Q7").Value = 2.6
$ws.Range("R7").Value = 1.48
$ws.Range("AD7").Value = 6.5
$ws.Range("AL7").Value = 41

# Row 8
$ws.Range("R8").Value = 1.57

# Row 9
$ws.Range("G9").Value = 1.7
$ws.Range("H9").Value = 4.2
$ws.Range("I9").Value = 4.33
$ws.Range("J9").Value = 2.2
$ws.Range("L9").Value = 4.5
$ws.Range("O9").Value = 1.14
$ws.Range("P9").Value = 5.5
$ws.Range("Q9").Value = 1.53
$ws.Range("R9").Value = 2.4
$ws.Range("S9").Value = 1.29
$ws.Range("T9").Value = 3.5
$ws.Range("U9").Value = 1.57
$ws.Range("V9").Value = 2.25
$ws.Range("W9").Value = 10
$ws.Range("X9").Value = 10
$ws.Range("AC9").Value = 17
$ws.Range("AD9").Value = 8.5
$ws.Range("AE9").Value = 13
$ws.Range("AF9").Value = 41
$ws.Range("AG9").Value = 126
$ws.Range("AJ9").Value = 15
$ws.Range("AK9").Value = 51
$ws.Range("AQ9").Value = 23
$ws.Range("AT9").Value = 3.5
$ws.Range("AU9").Value = 7.5
$ws.Range("BA9").Value = 81
$ws.Range("BD9").Value = 151

# Row 10
$ws.Range("G10").Value = 1.91
